$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 856.25
$ws.Range("I19").Value = 575
$ws.Range("K19").Value = 575
$ws.Range("M19").Value = -400

$ws.Range("H51").Value = 7999.5
$ws.Range("I51").Value = 7999.5
$ws.Range("K51").Value = 7999.5
$ws.Range("M51").Value = -7515.5

$ws.Range("H86").Value = 5000
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 5000
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 25000
$ws.Range("N89").Value = -36232

$ws.Range("H101").Value = 737.5
$ws.Range("I101").Value = 500
$ws.Range("J101").Value = 975
$ws.Range("K101").Value = 1500
$ws.Range("L101").Value = 2925
$ws.Range("M101").Value = 122
$ws.Range("N101").Value = -6169

$ws.Range("H106").Value = 2889.4443
$ws.Range("I106").Value = 2502.5
$ws.Range("K106").Value = 2502.5
$ws.Range("M106").Value = -1871.5

$ws.Range("H137").Value = 3309
$ws.Range("I137").Value = 2898.6667
$ws.Range("K137").Value = 8696.000100000001
$ws.Range("M137").Value = -6146.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 15002250
$ws.Range("I11").Value = 15002250
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 15002250
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -15002106
$ws.Range("N11").ClearContents()

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H97").Value = 1762.3077
$ws.Range("I97").Value = 1219.091
$ws.Range("K97").Value = 1219.091
$ws.Range("M97").Value = -723.0909999999999

$ws.Range("H101").Value = 24999
$ws.Range("J101").Value = 24999
$ws.Range("L101").Value = 24999
$ws.Range("N101").Value = -31489

$ws.Range("H102").Value = 1825.6666
$ws.Range("J102").Value = 1499.875
$ws.Range("L102").Value = 1499.875
$ws.Range("N102").Value = -4743.875

$ws.Range("H122").Value = 1613
$ws.Range("I122").Value = 1521.6875
$ws.Range("K122").Value = 4565.0625
$ws.Range("M122").Value = -2115.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1249.5
$ws.Range("I86").Value = 1249.5
$ws.Range("K86").Value = 1249.5
$ws.Range("M86").Value = -126.5

$ws.Range("H89").Value = 1249.5
$ws.Range("I89").Value = 1249.5
$ws.Range("K89").Value = 6247.5
$ws.Range("M89").Value = -631.5

$ws.Range("H95").Value = 15741.333
$ws.Range("J95").Value = 15741.333
$ws.Range("L95").Value = 15741.333
$ws.Range("N95").Value = -21233.333

$ws.Range("H132").Value = 99995
$ws.Range("J132").Value = 99995
$ws.Range("L132").Value = 99995
$ws.Range("N132").Value = -110115

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4384.5557
$ws.Range("I31").Value = 2096.4
$ws.Range("K31").Value = 2096.4
$ws.Range("M31").Value = -1801.4

$ws.Range("H34").Value = 4384.5557
$ws.Range("I34").Value = 2096.4
$ws.Range("K34").Value = 2096.4
$ws.Range("M34").Value = -1894.4

$ws.Range("H58").Value = 2679.4546
$ws.Range("I58").Value = 2871.875
$ws.Range("K58").Value = 2871.875
$ws.Range("M58").Value = -2668.875

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H105").Value = 1824.75
$ws.Range("I105").Value = 933
$ws.Range("K105").Value = 933
$ws.Range("M105").Value = 814

$ws.Range("H122").Value = 2748.3333
$ws.Range("I122").Value = 2175
$ws.Range("J122").Value = 3895
$ws.Range("K122").Value = 6525
$ws.Range("L122").Value = 11685
$ws.Range("M122").Value = -4075
$ws.Range("N122").Value = -16585

$ws.Range("H134").Value = 919.5
$ws.Range("I134").Value = 919.5
$ws.Range("K134").Value = 2758.5
$ws.Range("M134").Value = -223.5

$ws.Range("H136").Value = 2679.4546
$ws.Range("I136").Value = 2871.875
$ws.Range("K136").Value = 8615.625
$ws.Range("M136").Value = -6065.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 500535.5
$ws.Range("I4").Value = 500535.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1501606.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1501494.5
$ws.Range("N4").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H54").Value = 2000
$ws.Range("J54").Value = 2000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7118

$ws.Range("H123").Value = 600
$ws.Range("I123").Value = 600
$ws.Range("K123").Value = 1800
$ws.Range("M123").Value = 650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 49300.25
$ws.Range("J95").Value = 49300.25
$ws.Range("L95").Value = 49300.25
$ws.Range("N95").Value = -54792.25

$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815

$ws.Range("H126").Value = 1300
$ws.Range("I126").Value = 1300
$ws.Range("K126").Value = 3900
$ws.Range("M126").Value = -1430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1233.36
$ws.Range("I22").Value = 1286.75
$ws.Range("J22").Value = 1019.8
$ws.Range("K22").Value = 1286.75
$ws.Range("L22").Value = 1019.8
$ws.Range("M22").Value = -991.75
$ws.Range("N22").Value = -1609.8

$ws.Range("H27").Value = 1233.36
$ws.Range("I27").Value = 1286.75
$ws.Range("J27").Value = 1019.8
$ws.Range("K27").Value = 1286.75
$ws.Range("L27").Value = 1019.8
$ws.Range("M27").Value = -1179.75
$ws.Range("N27").Value = -1233.8

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H46").Value = 4548.021
$ws.Range("I46").Value = 2858.2
$ws.Range("K46").Value = 2858.2
$ws.Range("M46").Value = -2670.2

$ws.Range("H48").Value = 20000
$ws.Range("I48").Value = 20000
$ws.Range("K48").Value = 20000
$ws.Range("M48").Value = -19339

$ws.Range("H68").Value = 3249.8
$ws.Range("I68").Value = 2777.5557
$ws.Range("J68").Value = 7500
$ws.Range("K68").Value = 2777.5557
$ws.Range("L68").Value = 7500
$ws.Range("M68").Value = -2028.5557
$ws.Range("N68").Value = -8998

$ws.Range("H71").Value = 3249.8
$ws.Range("I71").Value = 2777.5557
$ws.Range("J71").Value = 7500
$ws.Range("K71").Value = 13887.7785
$ws.Range("L71").Value = 37500
$ws.Range("M71").Value = -10143.7785
$ws.Range("N71").Value = -44988

$ws.Range("H82").Value = 1724
$ws.Range("I82").Value = 1682.5555
$ws.Range("K82").Value = 1682.5555
$ws.Range("M82").Value = -1321.5555

$ws.Range("H85").Value = 1724
$ws.Range("I85").Value = 1682.5555
$ws.Range("K85").Value = 1682.5555
$ws.Range("M85").Value = -434.5554999999999

$ws.Range("H122").Value = 5537.4
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 4165.5557
$ws.Range("I136").Value = 4086.5715
$ws.Range("K136").Value = 12259.7145
$ws.Range("M136").Value = -9709.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 2501.5
$ws.Range("I96").Value = 3003
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 3003
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -1630
$ws.Range("N96").Value = -4746

$ws.Range("H100").Value = 506.42856
$ws.Range("I100").Value = 533.3333
$ws.Range("J100").Value = 345
$ws.Range("K100").Value = 1066.6666
$ws.Range("L100").Value = 690
$ws.Range("M100").Value = -525.6666
$ws.Range("N100").Value = -1772

$ws.Range("H101").Value = 17963.334
$ws.Range("J101").Value = 17963.334
$ws.Range("L101").Value = 17963.334
$ws.Range("N101").Value = -24453.334

$ws.Range("H111").Value = 40644
$ws.Range("J111").Value = 40644
$ws.Range("L111").Value = 40644
$ws.Range("N111").Value = -48824
